$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: Objetivos value replaced ---
$ws.Range("B10").Value = "5840730 - Antonio Jefferson da Silva Machado"
$ws.Range("C10").Value = "5840730 - Antonio Jefferson da Silva Machado"

# --- Row 13: now "Programa resumido:" / "Semestral" (gains an A-label) ---
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows.Item(13).RowHeight = 60

# --- Row 14: now just "Short syllabus:" label, values cleared ---
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14:C14").Clear()
$ws.Rows.Item(14).RowHeight = 60

# --- Row 15: now "Programa:" / "01/01/2012" (copy to avoid date coercion) ---
$ws.Range("A15").Value = "Programa:"
$ws.Range("B8:C8").Copy($ws.Range("B15:C15"))
$ws.Rows.Item(15).RowHeight = 120

# --- Row 16: now just "Syllabus:" label, values cleared ---
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16:C16").Clear()
$ws.Rows.Item(16).RowHeight = 120

# --- Row 17: now just "Avaliação:" label, default height ---
$ws.Range("A17").Value = "Avaliação:"
$ws.Rows.Item(17).EntireRow.AutoFit()

# --- Row 18: now "Método:" / docente name ---
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "5840730 - Antonio Jefferson da Silva Machado"
$ws.Range("C18").Value = "5840730 - Antonio Jefferson da Silva Machado"
$ws.Rows.Item(18).RowHeight = 60

# --- Row 19: label becomes "Critério:" ---
$ws.Range("A19").Value = "Critério:"

# --- Row 20: label becomes "Norma de recuperação:" ---
$ws.Range("A20").Value = "Norma de recuperação:"

# --- Row 21: label becomes "Bibliografia:", and height grows to 120 ---
$ws.Range("A21").Value = "Bibliografia:"
$ws.Rows.Item(21).RowHeight = 120

# --- Row 22 no longer exists: delete it entirely ---
$ws.Rows.Item(22).Delete()
